# Workbook: "memory map.xlsx"
# Adds a small "BIOS memory map" reference table to Arkusz2 (sheet 2),
# makes Arkusz2 the active/selected sheet, and moves the selection on
# Arkusz1 from F3 to A2:F2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- New reference table on Arkusz2 ---------------------------------
$ws2.Range("A1").Value = "początek biosu w pamięci"
$ws2.Range("B1").Value = "start PC"
$ws2.Range("C1").Value = "adres w pliku BIOS"
$ws2.Range("D1").Value = "przesunięcie"

$ws2.Range("A2").Value = "0xFE000"
$ws2.Range("B2").Value = "0xFFFF0"
$ws2.Range("C2").Value = "0x1FF0"
$ws2.Range("D2").Value = "0x1FF0"

# Column widths to fit the new content.
$ws2.Columns.Item(1).ColumnWidth = 29.82
$ws2.Columns.Item(2).ColumnWidth = 9.14
$ws2.Columns.Item(3).ColumnWidth = 18.14
$ws2.Columns.Item(4).ColumnWidth = 17.62

# Match the page setup used on Arkusz1.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Selection / active sheet ---------------------------------------
# Arkusz1's selection moves off of F3 onto the first data row.
$ws1.Range("A2:F2").Select()

# Arkusz2 becomes the active (selected) tab, with A2 selected.
$ws2.Activate()
$ws2.Range("A2").Select()
